$wb = $excel.ActiveWorkbook

# ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 3818.182
$ws.Range("I29").Value = 3937.5
$ws.Range("J29").Value = 3500
$ws.Range("K29").Value = 11812.5
$ws.Range("L29").Value = 10500
$ws.Range("M29").Value = -11531.5
$ws.Range("N29").Value = -11062
# Row 141
$ws.Range("H141").Value = 2676.65
$ws.Range("I141").Value = 2131.353
$ws.Range("K141").Value = 6394.059
$ws.Range("M141").Value = -1214.059

# ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 2871
$ws.Range("I63").Value = 2166.3333
$ws.Range("J63").Value = 7099
$ws.Range("K63").Value = 2166.3333
$ws.Range("L63").Value = 7099
$ws.Range("M63").Value = -1480.3333
$ws.Range("N63").Value = -8471
# Row 66
$ws.Range("H66").Value = 2871
$ws.Range("I66").Value = 2166.3333
$ws.Range("J66").Value = 7099
$ws.Range("K66").Value = 10831.6665
$ws.Range("L66").Value = 35495
$ws.Range("M66").Value = -7399.666499999999
$ws.Range("N66").Value = -42359

# BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 2366.7144
$ws.Range("I22").Value = 267
$ws.Range("J22").Value = 5166.3335
$ws.Range("K22").Value = 267
$ws.Range("L22").Value = 5166.3335
$ws.Range("M22").Value = -94
$ws.Range("N22").Value = -5512.3335
# Row 40
$ws.Range("H40").Value = 54500
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 54500
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 54500
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -55030
# Row 94
$ws.Range("H94").Value = 3531.037
$ws.Range("I94").Value = 3621.7273
$ws.Range("J94").Value = 3132
$ws.Range("K94").Value = 3621.7273
$ws.Range("L94").Value = 3132
$ws.Range("M94").Value = -3170.7273
$ws.Range("N94").Value = -4034

# CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 3428.8
$ws.Range("I132").Value = 3198.2
$ws.Range("J132").Value = 3659.4
$ws.Range("K132").Value = 9594.599999999999
$ws.Range("L132").Value = 10978.2
$ws.Range("M132").Value = -7064.599999999999
$ws.Range("N132").Value = -16038.2
# Row 134
$ws.Range("H134").Value = 4691.5
$ws.Range("I134").Value = 3538.75
$ws.Range("K134").Value = 10616.25
$ws.Range("M134").Value = -8081.25

# CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 2488
$ws.Range("J34").Value = 2720.4
$ws.Range("L34").Value = 8161.200000000001
$ws.Range("N34").Value = -8329.200000000001
# Row 44
$ws.Range("H44").Value = 7350
$ws.Range("J44").Value = 7350
$ws.Range("L44").Value = 22050
$ws.Range("N44").Value = -22846
# Row 113
$ws.Range("H113").Value = 898.1667
$ws.Range("I113").Value = 442
$ws.Range("J113").Value = 989.4
$ws.Range("K113").Value = 1326
$ws.Range("L113").Value = 2968.2
$ws.Range("M113").Value = 844
$ws.Range("N113").Value = -7308.2
# Row 117
$ws.Range("H117").Value = 1009
$ws.Range("I117").Value = 1009
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 3027
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 415
$ws.Range("N117").ClearContents()
# Row 118
$ws.Range("H118").Value = 1286
$ws.Range("I118").Value = 1181.3334
$ws.Range("K118").Value = 3544.0002
$ws.Range("M118").Value = -2301.0002
# Row 122
$ws.Range("H122").Value = 769624.7
$ws.Range("J122").Value = 1686.4667
$ws.Range("L122").Value = 15178.2003
$ws.Range("N122").Value = -20078.2003
# Row 123
$ws.Range("H123").Value = 13142.714
$ws.Range("I123").Value = 8499.5
$ws.Range("J123").Value = 15000
$ws.Range("K123").Value = 25498.5
$ws.Range("L123").Value = 45000
$ws.Range("M123").Value = -23048.5
$ws.Range("N123").Value = -49900
# Row 124
$ws.Range("H124").Value = 12291.091
$ws.Range("I124").Value = 10033.667
$ws.Range("J124").Value = 15000
$ws.Range("K124").Value = 30101.001
$ws.Range("L124").Value = 45000
$ws.Range("M124").Value = -25191.001
$ws.Range("N124").Value = -54820
# Row 125
$ws.Range("H125").Value = 2440.3333
$ws.Range("I125").Value = 2440.3333
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 7320.999899999999
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -2400.999899999999
$ws.Range("N125").ClearContents()
# Row 129
$ws.Range("H129").Value = 563165.4399999999
$ws.Range("I129").Value = 1253243.2
$ws.Range("J129").Value = 11103.2
$ws.Range("K129").Value = 3759729.6
$ws.Range("L129").Value = 33309.60000000001
$ws.Range("M129").Value = -3754729.6
$ws.Range("N129").Value = -43309.60000000001
# Row 130
$ws.Range("H130").Value = 9558.799999999999
$ws.Range("J130").Value = 15000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
# Row 131
$ws.Range("H131").Value = 6558.125
$ws.Range("J131").Value = 7314.5293
$ws.Range("L131").Value = 21943.5879
$ws.Range("N131").Value = -32023.5879

# GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 46
$ws.Range("H46").Value = 15000
$ws.Range("J46").Value = 15000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15312
# Row 57
$ws.Range("H57").Value = 21599.4
$ws.Range("J57").Value = 24499.25
$ws.Range("L57").Value = 24499.25
$ws.Range("N57").Value = -26139.25
# Row 64
$ws.Range("H64").Value = 25000
$ws.Range("I64").Value = 25000
$ws.Range("K64").Value = 25000
$ws.Range("M64").Value = -24752
# Row 67
$ws.Range("H67").Value = 25000
$ws.Range("I67").Value = 25000
$ws.Range("K67").Value = 25000
$ws.Range("M67").Value = -24142
# Row 80
$ws.Range("H80").Value = 2849.1428
$ws.Range("J80").Value = 3413.3333
$ws.Range("L80").Value = 3413.3333
$ws.Range("N80").Value = -5409.3333
# Row 83
$ws.Range("H83").Value = 2849.1428
$ws.Range("J83").Value = 3413.3333
$ws.Range("L83").Value = 17066.6665
$ws.Range("N83").Value = -27050.6665
# Row 102
$ws.Range("H102").Value = 14516.667
$ws.Range("J102").Value = 14000
$ws.Range("L102").Value = 14000
$ws.Range("N102").Value = -17244

# LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3089.2
$ws.Range("J22").Value = 3699.2
$ws.Range("L22").Value = 3699.2
$ws.Range("N22").Value = -4289.2
# Row 27
$ws.Range("H27").Value = 3089.2
$ws.Range("J27").Value = 3699.2
$ws.Range("L27").Value = 3699.2
$ws.Range("N27").Value = -3913.2
# Row 68
$ws.Range("H68").Value = 5452.2104
$ws.Range("I68").Value = 4253.879
$ws.Range("K68").Value = 4253.879
$ws.Range("M68").Value = -3504.879
# Row 71
$ws.Range("H71").Value = 5452.2104
$ws.Range("I71").Value = 4253.879
$ws.Range("K71").Value = 21269.395
$ws.Range("M71").Value = -17525.395
# Row 82
$ws.Range("H82").Value = 1748.8889
$ws.Range("I82").Value = 1592.5
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1592.5
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1231.5
$ws.Range("N82").Value = -3722
# Row 85
$ws.Range("H85").Value = 1748.8889
$ws.Range("I85").Value = 1592.5
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1592.5
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -344.5
$ws.Range("N85").Value = -5496
# Row 92
$ws.Range("H92").Value = 49900
$ws.Range("J92").Value = 49900
$ws.Range("L92").Value = 49900
$ws.Range("N92").Value = -54892
# Row 122
$ws.Range("H122").Value = 6362
$ws.Range("I122").Value = 3864.4
$ws.Range("K122").Value = 11593.2
$ws.Range("M122").Value = -9143.200000000001
